# Add "Sheet2" after "Sheet1" and populate it with intent/queries/responses
# data (TAT + PLM Defects Q&A pairs), matching the author's data-update commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: selection moves to A1:C1 (no longer the active tab) ---
$ws1.Activate()
$ws1.Range("A1:C1").Select()

# --- Create Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$rows = @(
    @('intent', 'queries', 'responses'),
    @($null, 'What exactly is TAT?', 'The Turnaround Time (TAT) is calculated as the duration between the Resolve Date and the Open Date, representing the time taken to address and resolve defects.'),
    @($null, 'Could you please explain what TAT means?', 'The Turnaround Time (TAT) measures the elapsed time from the moment a defect is opened until it is resolved.'),
    @($null, 'Can you provide a definition for TAT?', 'TAT refers to the time it takes to fix and resolve defects, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'I''m not familiar with TAT, could you provide some information?', 'The time required to resolve defects, known as Turnaround Time (TAT), is determined by subtracting the Open Date from the Resolve Date.'),
    @($null, 'I''m curious about TAT, can you tell me more?', 'When we talk about Turnaround Time (TAT), we''re referring to the timeframe it takes to address and resolve defects, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'What does TAT refer to?', 'The duration it takes to fix defects, commonly referred to as Turnaround Time (TAT), is obtained by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Can you give me some insights into what TAT is?', 'Turnaround Time (TAT) is a metric used to quantify the time it takes to resolve defects, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'I''d like to know more about TAT.', 'The amount of time taken to resolve defects, known as Turnaround Time (TAT), is determined by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Could you elaborate on what TAT entails?', 'Turnaround Time (TAT) represents the time required to fix defects, calculated by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Can you shed some light on TAT?', 'The duration between the Resolve Date and the Open Date is referred to as Turnaround Time (TAT), which signifies the time taken to resolve defects.'),
    @($null, 'I''m interested in learning about TAT.', 'In the context of defect resolution, Turnaround Time (TAT) is the period it takes to address and fix issues, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'What is the significance of TAT?', 'The time taken to resolve defects, denoted as Turnaround Time (TAT), is calculated by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Can you give me a brief overview of TAT?', 'TAT, or Turnaround Time, quantifies the duration between the Resolve Date and the Open Date, indicating how long it took to resolve defects.'),
    @($null, 'I''m looking for information about TAT, can you help?', 'When we talk about Turnaround Time (TAT), we''re referring to the timeframe it takes to address and fix defects, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'What can you tell me about TAT?', 'Turnaround Time (TAT) measures the elapsed time from when a defect is opened until it is resolved, indicating how long it took to fix the issues.'),
    @($null, 'I''m curious to know the meaning of TAT.', 'The duration it takes to resolve defects, often referred to as Turnaround Time (TAT), is obtained by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Can you provide some context for TAT?', 'Turnaround Time (TAT) represents the time taken to address and resolve defects, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'What purpose does TAT serve?', 'TAT, or Turnaround Time, is a metric that captures the time it takes to resolve defects, calculated by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Please provide a description of TAT.', 'The amount of time taken to resolve defects, known as Turnaround Time (TAT), is determined by subtracting the Open Date from the Resolve Date.'),
    @($null, 'I''d appreciate it if you could explain TAT to me.', 'Turnaround Time (TAT) is the duration between the Resolve Date and the Open Date, indicating the time taken to resolve defects.'),
    @($null, 'Could you give me a synopsis of TAT?', 'The time taken to fix and resolve defects, known as Turnaround Time (TAT), is calculated by subtracting the Open Date from the Resolve Date.'),
    @($null, 'Can you offer some details about TAT?', 'Turnaround Time (TAT) measures the efficiency of defect resolution by calculating the duration between the Resolve Date and the Open Date.'),
    @($null, 'I''d like to gain a better understanding of TAT.', 'TAT, or Turnaround Time, is an important metric that helps gauge the speed and effectiveness of defect resolution, calculated as the difference between the Resolve Date and the Open Date.'),
    @($null, 'What are the key aspects of TAT?', 'The Turnaround Time (TAT) metric reflects the responsiveness of defect resolution, indicating the time it took to address and resolve issues.'),
    @($null, 'Could you break down the concept of TAT for me?', 'When we refer to Turnaround Time (TAT), we''re focusing on the timeframe it took to resolve defects, derived from the difference between the Resolve Date and the Open Date.'),
    @($null, 'I''m seeking clarification on what TAT entails.', 'Turnaround Time (TAT) measures the elapsed time from when a defect is opened until it is resolved, indicating how long it took to fix the issues.'),
    @($null, '1. Can you explain the concept of PLM defects?', '1. PLM Defects refer to the flaws or issues identified within a Product Lifecycle Management (PLM) system.'),
    @($null, '2. Could you provide some insights into PLM defects?', '2. PLM Defects are defects or problems that occur within the context of Product Lifecycle Management (PLM), affecting the system''s functionality or performance.'),
    @($null, '3. How would you define PLM defects?', '3. PLM Defects are the defects or anomalies found in the Product Lifecycle Management (PLM) process, which may impact the system''s efficiency or accuracy.'),
    @($null, '4. Can you elaborate on the notion of PLM defects?', '4. When we talk about PLM Defects, we are referring to the issues or errors encountered within the Product Lifecycle Management (PLM) framework, affecting its overall operation.'),
    @($null, '5. What do we mean by PLM defects?', '5. PLM Defects pertain to the defects, bugs, or shortcomings discovered within a Product Lifecycle Management (PLM) system, potentially hindering its effectiveness or smooth functioning.'),
    @($null, '6. Can you shed some light on PLM defects?', '1. PLM Defects are deviations or discrepancies identified during the Product Lifecycle Management (PLM) process, indicating areas that require attention or improvement.'),
    @($null, '7. What constitutes PLM defects?', '2. When discussing PLM Defects, we are referring to the instances where the Product Lifecycle Management (PLM) system fails to meet the expected standards or specifications.'),
    @($null, '8. Can you give me an overview of PLM defects?', '3. PLM Defects encompass the issues, faults, or irregularities detected within the Product Lifecycle Management (PLM) framework, highlighting the need for corrective actions.'),
    @($null, '9. In the context of PLM, what are defects?', '4. PLM Defects represent the non-conformances or non-compliances encountered within the Product Lifecycle Management (PLM) domain, necessitating thorough investigation and resolution.'),
    @($null, $null, '5. The term PLM Defects refers to the shortcomings, flaws, or deficiencies observed in the Product Lifecycle Management (PLM) system, requiring remedial measures.'),
    @($null, '11. How are PLM defects defined and categorized?', '6. PLM Defects indicate the instances where the Product Lifecycle Management (PLM) system deviates from the desired state, requiring troubleshooting and rectification.'),
    @($null, '12. What is the significance of PLM defects?', '7. PLM Defects are the documented instances of deviations or malfunctions encountered within the Product Lifecycle Management (PLM) processes, aiming for continuous improvement.'),
    @($null, $null, '8. When we talk about PLM Defects, we are referring to the anomalies or irregularities identified within the Product Lifecycle Management (PLM) system, demanding corrective actions.'),
    @($null, $null, '9. PLM Defects signify the occurrences where the Product Lifecycle Management (PLM) system fails to adhere to the predefined quality standards, necessitating investigation and resolution.'),
    @($null, $null, '10. PLM Defects encompass the recorded instances of non-conforming or erroneous behaviors within the Product Lifecycle Management (PLM) operations, requiring analysis and rectification.'),
    @($null, $null, '11. When addressing PLM Defects, we are referring to the imperfections, errors, or glitches detected in the Product Lifecycle Management (PLM) framework, warranting attention and remediation.'),
    @($null, $null, '12. PLM Defects represent the documented observations of discrepancies or deviations encountered within the Product Lifecycle Management (PLM) processes, aiming for process optimization.'),
    @($null, $null, '13. PLM Defects highlight the areas where the Product Lifecycle Management (PLM) system falls short of meeting the defined quality criteria, necessitating corrective measures.'),
    @($null, $null, '14. When discussing PLM Defects, we are referring to the instances of non-compliance, malfunctions, or suboptimal performance within the Product Lifecycle Management (PLM) domain.'),
    @($null, $null, '15. PLM Defects encompass the identified issues, errors, or inconsistencies observed during the Product Lifecycle Management (PLM) activities, aiming for continuous quality enhancement.')
)

$r = 1
foreach ($row in $rows) {
    if ($row[0]) { $ws2.Cells.Item($r, 1).Value = $row[0] }
    if ($row[1]) { $ws2.Cells.Item($r, 2).Value = $row[1] }
    if ($row[2]) { $ws2.Cells.Item($r, 3).Value = $row[2] }
    $r = $r + 1
}

# Autofit columns to match bestFit widths, like the original author's sheet
$ws2.Columns("A:C").AutoFit()

# Sheet2 becomes the active sheet/tab, with B47 selected (last data row)
$ws2.Activate()
$ws2.Range("B47").Select()
